# Updates from yesterday's status meeting
# Fill in the three previously-blank action items (rows 93-95) on the
# Action_Items sheet, then leave the selection where the user ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Action_Items")

# Dates in this workbook use the 1904 date system (workbookPr date1904="1"),
# and the meeting date is 2013-06-04, which is serial 39967 under that
# system. We poke the raw serial number in directly (the cell already has
# the date number format applied via its style) instead of assigning a
# date/string literal, since date literals get converted assuming the 1900
# system regardless of the workbook's date1904 setting.
$meetingDate = 39967

# --- Row 93 (Item #92) ---------------------------------------------------
$ws.Range("B93").Value = "Discuss the documentation update process"
$ws.Range("C93").Value = "Mike Hunter and Jill Hadfield"
$ws.Range("D93").Value = $meetingDate
$ws.Range("E93").Value = "Assigned"

# --- Row 94 (Item #93) ---------------------------------------------------
$ws.Range("B94").Value = "Find the past FISMA documentation and confirm the tracker issues"
$ws.Range("C94").Value = "Mike Hunter"
$ws.Range("D94").Value = $meetingDate
$ws.Range("E94").Value = "Assigned"

# --- Row 95 (Item #94) ---------------------------------------------------
$ws.Range("B95").Value = "Add a %FTE column to the project plans and populate it going forward"
$ws.Range("C95").Value = "Mike Hunter and Shine Jacob"
$ws.Range("D95").Value = $meetingDate
$ws.Range("E95").Value = "Assigned"

# --- Window / view state, matching where the author left the cursor -----
$ws.Activate()
$ws.Range("E96").Select()

# Scroll so column B is pinned at the left edge of the viewport (was column
# A) and nudge the main Excel window's on-screen position, matching the
# author's window placement when the file was saved.
try { $excel.ActiveWindow.ScrollColumn = 2 } catch {}
try { $excel.ActiveWindow.ScrollRow = 75 } catch {}
try { $excel.ActiveWindow.Left = 3780 } catch {}
